# Auto-generated Excel COM-interop script
# Applies numeric value updates (and a couple of cell insert/deletes)
# to the Ifrit_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as described by the source diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 6 (Leve Item ID 4564)
$ws.Range("H6").Value = 87.5
$ws.Range("I6").Value = 87.5
$ws.Range("K6").Value = 262.5
$ws.Range("M6").Value = -150.5

# ALC row 28 (Leve Item ID 27772)
$ws.Range("H28").Value = 1080.5883
$ws.Range("I28").Value = 1136.6666
$ws.Range("J28").Value = 946
$ws.Range("K28").Value = 1136.6666
$ws.Range("L28").Value = 946
$ws.Range("M28").Value = -651.6666
$ws.Range("N28").Value = -1916

# ALC row 76 (Leve Item ID 12602)
$ws.Range("H76").Value = 3106.4888
$ws.Range("I76").Value = 3106.4888
$ws.Range("K76").Value = 3106.4888
$ws.Range("M76").Value = -2791.4888

# ALC row 79 (Leve Item ID 12602)
$ws.Range("H79").Value = 3106.4888
$ws.Range("I79").Value = 3106.4888
$ws.Range("K79").Value = 3106.4888
$ws.Range("M79").Value = -2014.4888

# ALC row 86 (Leve Item ID 12603)
$ws.Range("H86").Value = 3185.5
$ws.Range("I86").Value = 2950
$ws.Range("J86").Value = 3421
$ws.Range("K86").Value = 2950
$ws.Range("L86").Value = 3421
$ws.Range("M86").Value = -1827
$ws.Range("N86").Value = -5667

# ALC row 89 (Leve Item ID 12603)
$ws.Range("H89").Value = 3185.5
$ws.Range("I89").Value = 2950
$ws.Range("J89").Value = 3421
$ws.Range("K89").Value = 14750
$ws.Range("L89").Value = 17105
$ws.Range("M89").Value = -9134
$ws.Range("N89").Value = -28337

# ALC row 107 (Leve Item ID 27766)
$ws.Range("H107").Value = 1154.8125
$ws.Range("I107").Value = 1154.8125
$ws.Range("K107").Value = 1154.8125
$ws.Range("M107").Value = 765.1875

# ALC row 113 (Leve Item ID 27775)
$ws.Range("H113").Value = 2088.25
$ws.Range("I113").Value = 2500
$ws.Range("J113").Value = 1951
$ws.Range("K113").Value = 2500
$ws.Range("L113").Value = 1951
$ws.Range("M113").Value = 754
$ws.Range("N113").Value = -8459

# ALC row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 5139
$ws.Range("I137").Value = 8318.866
$ws.Range("J137").Value = 2628.5789
$ws.Range("K137").Value = 24956.598
$ws.Range("L137").Value = 7885.736699999999
$ws.Range("M137").Value = -22406.598
$ws.Range("N137").Value = -12985.7367

# ALC row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 1155900.2
$ws.Range("I138").Value = 2550.3
$ws.Range("J138").Value = 1705114.5
$ws.Range("K138").Value = 7650.900000000001
$ws.Range("L138").Value = 5115343.5
$ws.Range("M138").Value = -2510.900000000001
$ws.Range("N138").Value = -5125623.5

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 4202159
$ws.Range("I2").Value = 537.5
$ws.Range("J2").Value = 9804321
$ws.Range("K2").Value = 537.5
$ws.Range("L2").Value = 9804321
$ws.Range("M2").Value = -424.5
$ws.Range("N2").Value = -9804547

# ARM row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 11763.884
$ws.Range("I32").Value = 5255.8184
$ws.Range("J32").Value = 29661.062
$ws.Range("K32").Value = 5255.8184
$ws.Range("L32").Value = 29661.062
$ws.Range("M32").Value = -4968.8184
$ws.Range("N32").Value = -30235.062

# ARM row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 45295.89
$ws.Range("I110").Value = 924.5454999999999
$ws.Range("K110").Value = 924.5454999999999
$ws.Range("M110").Value = 1120.4545

# ARM row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 4202159
$ws.Range("I116").Value = 537.5
$ws.Range("J116").Value = 9804321
$ws.Range("K116").Value = 537.5
$ws.Range("L116").Value = 9804321
$ws.Range("M116").Value = 1756.5
$ws.Range("N116").Value = -9808909

# ARM row 138 (Leve Item ID 42350)
$ws.Range("H138").Value = 42195.668
$ws.Range("I138").Value = 54787
$ws.Range("J138").Value = 35900
$ws.Range("K138").Value = 54787
$ws.Range("L138").Value = 35900
$ws.Range("N138").Value = -46180
$ws.Range("M138").Value = -49647

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 4202159
$ws.Range("I3").Value = 537.5
$ws.Range("J3").Value = 9804321
$ws.Range("K3").Value = 537.5
$ws.Range("L3").Value = 9804321
$ws.Range("M3").Value = -423.5
$ws.Range("N3").Value = -9804549

# BSM row 20 (Leve Item ID 14149)
$ws.Range("H20").Value = 1288.963
$ws.Range("I20").Value = 1330.5
$ws.Range("J20").Value = 1244.2307
$ws.Range("K20").Value = 1330.5
$ws.Range("L20").Value = 1244.2307
$ws.Range("M20").Value = -1083.5
$ws.Range("N20").Value = -1738.2307

# BSM row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 1786.5769
$ws.Range("I105").Value = 1786.7222
$ws.Range("J105").Value = 1786.25
$ws.Range("K105").Value = 1786.7222
$ws.Range("L105").Value = 1786.25
$ws.Range("M105").Value = -39.72219999999993
$ws.Range("N105").Value = -5280.25

# BSM row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 672692.3
$ws.Range("I107").Value = 830055.2
$ws.Range("K107").Value = 830055.2
$ws.Range("M107").Value = -828135.2

$ws = $wb.Worksheets.Item("CRP")
# CRP row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 11365048
$ws.Range("I122").Value = 17858100
$ws.Range("K122").Value = 53574300
$ws.Range("M122").Value = -53571850

$ws = $wb.Worksheets.Item("CUL")
# CUL row 68 (Leve Item ID 12895)
$ws.Range("H68").Value = 873.6437
$ws.Range("I68").Value = 663.2
$ws.Range("J68").Value = 1099.119
$ws.Range("K68").Value = 1989.6
$ws.Range("L68").Value = 3297.357
$ws.Range("M68").Value = -1178.6
$ws.Range("N68").Value = -4919.357

# CUL row 71 (Leve Item ID 12895)
$ws.Range("H71").Value = 873.6437
$ws.Range("I71").Value = 663.2
$ws.Range("J71").Value = 1099.119
$ws.Range("K71").Value = 5968.8
$ws.Range("L71").Value = 9892.071
$ws.Range("M71").Value = -1912.8
$ws.Range("N71").Value = -18004.071

# CUL row 96 (Leve Item ID 19816)
$ws.Range("H96").Value = 2551.5
$ws.Range("J96").Value = 2551.5
$ws.Range("L96").Value = 7654.5
$ws.Range("N96").Value = -11772.5

# CUL row 105 (Leve Item ID 19814)
$ws.Range("H105").Value = 227002240
$ws.Range("J105").Value = 227002240
$ws.Range("L105").Value = 681006720
$ws.Range("N105").Value = -681011962

# CUL row 110 (Leve Item ID 27857)
$ws.Range("H110").Value = 1013.5
$ws.Range("I110").Value = 1013.5
$ws.Range("K110").Value = 3040.5
$ws.Range("M110").Value = 1049.5

# CUL row 114 (Leve Item ID 27865)
$ws.Range("H114").Value = 1306.0454
$ws.Range("I114").Value = 300.27274
$ws.Range("J114").Value = 2311.818
$ws.Range("K114").Value = 900.81822
$ws.Range("L114").Value = 6935.454000000001
$ws.Range("M114").Value = 2353.18178
$ws.Range("N114").Value = -13443.454

# CUL row 118 (Leve Item ID 27872)
$ws.Range("H118").Value = 959.8
$ws.Range("I118").Value = 224.75
$ws.Range("K118").Value = 674.25
$ws.Range("M118").Value = 568.75

# CUL row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 1494400.4
$ws.Range("J131").Value = 1853896.5
$ws.Range("L131").Value = 5561689.5
$ws.Range("N131").Value = -5571769.5

# CUL row 140 (Leve Item ID 44097)
$ws.Range("H140").Value = 23685032
$ws.Range("I140").Value = 23685032
$ws.Range("K140").Value = 71055096
$ws.Range("M140").Value = -71049916

$ws = $wb.Worksheets.Item("GSM")
# GSM row 5 (Leve Item ID 1681)
$ws.Range("H5").Value = 2113.375
$ws.Range("I5").Value = 1900
$ws.Range("J5").Value = 2117.9148
$ws.Range("K5").Value = 1900
$ws.Range("L5").Value = 2117.9148
$ws.Range("N5").Value = -2341.9148
$ws.Range("M5").Value = -1788

# GSM row 9 (Leve Item ID 1683)
$ws.Range("H9").Value = 3000
$ws.Range("J9").Value = 3000
$ws.Range("L9").Value = 3000
$ws.Range("N9").Value = -2341.9148
$ws.Range("M9").Value = -1788

# GSM row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 1256.4286
$ws.Range("I113").Value = 1199
$ws.Range("J113").Value = 1333
$ws.Range("K113").Value = 1199
$ws.Range("L113").Value = 1333
$ws.Range("M113").Value = 971
$ws.Range("N113").Value = -5673

$ws = $wb.Worksheets.Item("LTW")
# LTW row 9 (Leve Item ID 1685)
$ws.Range("H9").Value = 487.5
$ws.Range("I9").Value = 487.5
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 487.5
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -263.5
$ws.Range("N9").ClearContents()

# LTW row 30 (Leve Item ID 1688)
$ws.Range("H30").Value = 4508
$ws.Range("I30").Value = 4508
$ws.Range("K30").Value = 4508
$ws.Range("M30").Value = -4400

# LTW row 61 (Leve Item ID 27740)
$ws.Range("H61").Value = 4568
$ws.Range("I61").Value = 1852
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 1852
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -1650
$ws.Range("N61").Value = -10404

# LTW row 82 (Leve Item ID 12565)
$ws.Range("H82").Value = 3052.2307
$ws.Range("I82").Value = 2954.1428
$ws.Range("J82").Value = 3166.6667
$ws.Range("K82").Value = 2954.1428
$ws.Range("L82").Value = 3166.6667
$ws.Range("M82").Value = -2593.1428
$ws.Range("N82").Value = -3888.6667

# LTW row 85 (Leve Item ID 12565)
$ws.Range("H85").Value = 3052.2307
$ws.Range("I85").Value = 2954.1428
$ws.Range("J85").Value = 3166.6667
$ws.Range("K85").Value = 2954.1428
$ws.Range("L85").Value = 3166.6667
$ws.Range("M85").Value = -1706.1428
$ws.Range("N85").Value = -5662.6667

# LTW row 113 (Leve Item ID 27740)
$ws.Range("H113").Value = 4568
$ws.Range("I113").Value = 1852
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 1852
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = 318
$ws.Range("N113").Value = -14340

# LTW row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 7135.853
$ws.Range("I132").Value = 9264.157999999999
$ws.Range("J132").Value = 4440
$ws.Range("K132").Value = 27792.474
$ws.Range("L132").Value = 13320
$ws.Range("M132").Value = -25262.474
$ws.Range("N132").Value = -18380

# LTW row 141 (Leve Item ID 42487)
$ws.Range("H141").Value = 48975
$ws.Range("J141").Value = 48975
$ws.Range("L141").Value = 48975
$ws.Range("N141").Value = -59335

$ws = $wb.Worksheets.Item("WVR")
# WVR row 75 (Leve Item ID 11957)
$ws.Range("H75").Value = 30750
$ws.Range("I75").Value = 20000
$ws.Range("J75").Value = 34333.332
$ws.Range("K75").Value = 20000
$ws.Range("L75").Value = 34333.332
$ws.Range("M75").Value = -19064
$ws.Range("N75").Value = -36205.332

# WVR row 78 (Leve Item ID 11957)
$ws.Range("H78").Value = 30750
$ws.Range("I78").Value = 20000
$ws.Range("J78").Value = 34333.332
$ws.Range("K78").Value = 60000
$ws.Range("L78").Value = 102999.996
$ws.Range("M78").Value = -55320
$ws.Range("N78").Value = -112359.996

# WVR row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 640.8333
$ws.Range("I113").Value = 600
$ws.Range("J113").Value = 654.44446
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 1963.33338
$ws.Range("M113").Value = 370
$ws.Range("N113").Value = -6303.33338
